# Updated legacy GSC export data:
# The "Chart" sheet had a stray/incomplete row for 2025-09-15 (row 2) whose
# Indexed/Not-indexed values were accidentally stored as text and whose
# Impressions value duplicated the following day's row. That row is removed
# and all the data below it shifts up by one row (the sheet now ends at
# row 69 instead of row 70). No other sheets need structural changes -
# their shared-string references simply follow the shift automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the obsolete 2025-09-15 row (row 2), shifting everything below it up.
$ws.Range("A2:D2").Delete()
